$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Force Price/Volume columns to be treated as text so numeric-looking
# strings (e.g. "0.116", "7.88") are not coerced into floating point
# numbers by Excel, then restore the original (default) formatting.
$priceVolRange = $ws.Range("D2:E51")
$priceVolRange.NumberFormat = "@"

$ws.Range('D2').Value = '65.384.05'
$ws.Range('E2').Value = '  -1.88%  '
$ws.Range('D3').Value = '3.407.26'
$ws.Range('E3').Value = '  -1.43%  '
$ws.Range('E4').Value = '  -0.05%  '
$ws.Range('D5').Value = '593.58'
$ws.Range('E5').Value = '  -1.46%  '
$ws.Range('D6').Value = '142.34'
$ws.Range('E6').Value = '  -3.48%  '
$ws.Range('B7').Value = 'USDC'
$ws.Range('C7').Value = 'https://coinranking.com/coin/aKzUVe4Hh_CON+usdc-usdc'
$ws.Range('D7').Value = '0.999'
$ws.Range('E7').Value = '  -0.13%  '
$ws.Range('B8').Value = 'LidoStakedEther'
$ws.Range('C8').Value = 'https://coinranking.com/coin/VINVMYf0u+lidostakedether-steth'
$ws.Range('D8').Value = '3.408.42'
$ws.Range('E8').Value = '  -1.32%  '
$ws.Range('E9').Value = '  -3.31%  '
$ws.Range('E10').Value = '  -4.78%  '
$ws.Range('D11').Value = '7.88'
$ws.Range('E11').Value = '  +5.61%  '
$ws.Range('E12').Value = '  -3.88%  '
$ws.Range('D13').Value = '3.990.68'
$ws.Range('E13').Value = '  -1.25%  '
$ws.Range('D14').Value = '0.0000200'
$ws.Range('E14').Value = '  -5.86%  '
$ws.Range('D15').Value = '29.86'
$ws.Range('E15').Value = '  -5.38%  '
$ws.Range('B16').Value = 'TRON'
$ws.Range('C16').Value = 'https://coinranking.com/coin/qUhEFk1I61atv+tron-trx'
$ws.Range('D16').Value = '0.116'
$ws.Range('E16').Value = '  -0.58%  '
$ws.Range('B17').Value = 'WrappedEther'
$ws.Range('C17').Value = 'https://coinranking.com/coin/Mtfb0obXVh59u+wrappedether-weth'
$ws.Range('D17').Value = '3.405.06'
$ws.Range('E17').Value = '  -1.40%  '
$ws.Range('B18').Value = 'WrappedBTC'
$ws.Range('C18').Value = 'https://coinranking.com/coin/x4WXHge-vvFY+wrappedbtc-wbtc'
$ws.Range('D18').Value = '65.466.76'
$ws.Range('E18').Value = '  -1.96%  '
$ws.Range('D19').Value = '10.37'
$ws.Range('E19').Value = '  +4.06%  '
$ws.Range('E20').Value = '  -4.41%  '
$ws.Range('E21').Value = '  -2.77%  '
$ws.Range('D22').Value = '417.12'
$ws.Range('E22').Value = '  -4.86%  '
$ws.Range('D23').Value = '0.581'
$ws.Range('E23').Value = '  -4.74%  '
$ws.Range('D24').Value = '77.09'
$ws.Range('E24').Value = '  -1.19%  '
$ws.Range('D26').Value = '3.542.89'
$ws.Range('E26').Value = '  -1.45%  '
$ws.Range('D27').Value = '0.0000111'
$ws.Range('E27').Value = '  -7.46%  '
$ws.Range('D28').Value = '9.25'
$ws.Range('E28').Value = '  -5.92%  '
$ws.Range('D29').Value = '7.82'
$ws.Range('E29').Value = '  -6.69%  '
$ws.Range('D30').Value = '2.42'
$ws.Range('E30').Value = '  -2.21%  '
$ws.Range('E31').Value = '  +0.07%  '
$ws.Range('D32').Value = '0.162'
$ws.Range('E32').Value = '  -2.15%  '
$ws.Range('E33').Value = '  -8.10%  '
$ws.Range('D34').Value = '24.64'
$ws.Range('E34').Value = '  -2.72%  '
$ws.Range('E36').Value = '  -4.82%  '
$ws.Range('D37').Value = '5.55'
$ws.Range('E37').Value = '  -8.98%  '
$ws.Range('D38').Value = '7.60'
$ws.Range('E38').Value = '  -3.67%  '
$ws.Range('E39').Value = '  +0.06%  '
$ws.Range('D40').Value = '173.59'
$ws.Range('E40').Value = '  +0.25%  '
$ws.Range('D41').Value = '0.0864'
$ws.Range('E41').Value = '  -2.27%  '
$ws.Range('D42').Value = '5.06'
$ws.Range('E42').Value = '  -5.49%  '
$ws.Range('D43').Value = '0.868'
$ws.Range('E43').Value = '  -1.34%  '
$ws.Range('E44').Value = '  -11.25%  '
$ws.Range('D45').Value = '45.54'
$ws.Range('E45').Value = '  -1.14%  '
$ws.Range('D46').Value = '26.76'
$ws.Range('E46').Value = '  -7.12%  '
$ws.Range('E47').Value = '  -5.08%  '
$ws.Range('D48').Value = '7.09'
$ws.Range('E48').Value = '  -4.95%  '
$ws.Range('E49').Value = '  -6.71%  '
$ws.Range('D50').Value = '0.921'
$ws.Range('E50').Value = '  -6.51%  '
$ws.Range('E51').Value = '  -4.46%  '

$priceVolRange.ClearFormats()

